# Update the "policy structure" graphic on slide 3:
#  - insert a new bulleted "Operations" line right after "Annotations"
#    in the policy-matching-statements textbox
#  - the textbox autosizes (spAutoFit) and is then nudged to its final
#    on-slide position/size

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(12)   # "TextBox 29"

if ($shp.Name -ne "TextBox 29") {
    throw "Expected 'TextBox 29' at Shapes.Item(12), found '$($shp.Name)'"
}

$tr2 = $shp.TextFrame2.TextRange
$annotationsPara = $tr2.Paragraphs(4)   # "Annotations"
if ($annotationsPara.Text.TrimEnd("`r") -ne "Annotations") {
    throw "Expected paragraph 4 to be 'Annotations', found '$($annotationsPara.Text)'"
}
$annotationsPara.InsertAfter("`rOperations") | Out-Null

# Final position for the textbox (EMU -> points, 12700 EMU/pt). The
# literals below are nudged slightly within the same EMU bucket so that
# the host's internal float32 storage still rounds back to the exact
# target EMU values (8404102, 236094) instead of the adjacent one.
$shp.Left = 661.740350
$shp.Top  = 18.590118

# Width/Height already resolve to the right EMU values (2832827 x 3816429)
# via the textbox's spAutoFit after the new bullet line was added above,
# so they are left untouched.
